$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81-85 down to 82-86
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with data (mirrors the record previously at row 81,
# with updated Fecha/Volumen/Precios as per the commit)
$ws.Range("A81").Value = 9
$ws.Range("B81").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C81").Value = "Metropolitana"
$ws.Range("D81").Value = 44615
$ws.Range("E81").Value = 13
$ws.Range("F81").Value = 100112005
$ws.Range("G81").Value = "Puerro"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 79
$ws.Range("K81").Value = 7000
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = 7000
$ws.Range("N81").Value = "$/paquete 20 unidades"
$ws.Range("O81").Value = "Provincia de Chacabuco"
$ws.Range("P81").Value = 350
$ws.Range("Q81").Value = 20
$ws.Range("R81").Value = "Hortaliza"
